$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Increment "dias_defasagem" (column E) by 1 for each data row (2..41),
# reflecting one more day elapsed since the stacked bar chart data was refreshed.
for ($r = 2; $r -le 41; $r++) {
    $cell = $ws.Cells.Item($r, 5)  # Column E
    $cell.Value2 = $cell.Value2 + 1
}
